$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (strike counts) for rows 2-80, replacing the previous
# Strike#-derived values in column G after regenerating save_data with
# std/mean recalculated and s_vals written.
$kValues = @(1,0,0,0,2,2,3,2,1,1,2,0,1,1,1,0,2,1,1,0,2,0,1,1,1,1,0,1,2,0,1,0,0,0,0,0,1,1,0,2,3,1,0,1,0,0,1,0,0,1,2,0,2,2,0,1,1,1,1,1,1,1,1,1,1,2,0,3,0,1,0,3,2,0,1,3,3,1,1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
